# Updated week 10 wins.
# Fill in the Week 10 (column N) win totals on the "NFL" sheet for the
# teams/weeks that were missing them. The sheet's A-column CONCATENATE
# formulas and the summary rows (36-39, SUMIF totals per person) recalc
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NFL")

$ws.Range("N3").Value  = 1   # Greg / 49ers
$ws.Range("N5").Value  = 3   # Tim  / Bengals
$ws.Range("N10").Value = 4   # Jeff / Cardinals
$ws.Range("N11").Value = 4   # Jeff / Chargers
$ws.Range("N14").Value = 8   # Greg / Cowboys
$ws.Range("N15").Value = 5   # Tim  / Dolphins
$ws.Range("N18").Value = 6   # Jeff / Giants
$ws.Range("N24").Value = 7   # Greg / Patriots
$ws.Range("N30").Value = 6   # Zach / Seahawks
$ws.Range("N31").Value = 4   # Tim  / Steelers

# Leave the view the way the author left it: NFL sheet focused on N24,
# NBA sheet focused on F25.
$ws.Activate()
$ws.Range("N24").Select()

$wsNba = $wb.Worksheets.Item("NBA")
$wsNba.Activate()
$wsNba.Range("F25").Select()

$ws.Activate()
